$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grading for Q14 "whoPurchasedProduct() method" (row 22): 5 points deducted,
# with a comment explaining why.
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = " -5 for wrong logic"

# Grading for Q15 "getProducts()" (row 23): leave points as-is, just add a
# comment.
$ws.Range("F23").Value = " -6 for incomplete method"

# Move the viewport / active selection to where the grader was last working.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("F23").Select()
